$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Cells changing from numeric to text ("0" / "***.*") ---
$ws.Cells.Item(14, 4).Copy($ws.Cells.Item(15, 4))
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(15, 5))
$ws.Cells.Item(14, 4).Copy($ws.Cells.Item(26, 4))
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(26, 5))

# --- Cells changing from text to numeric (copy style from a same-style numeric cell, then set value) ---
$ws.Cells.Item(14, 7).Copy($ws.Cells.Item(23, 4))
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(14, 8).Copy($ws.Cells.Item(23, 5))
$ws.Cells.Item(23, 5).Value = -100
$ws.Cells.Item(14, 7).Copy($ws.Cells.Item(23, 10))
$ws.Cells.Item(23, 10).Value = 2
$ws.Cells.Item(14, 8).Copy($ws.Cells.Item(23, 11))
$ws.Cells.Item(23, 11).Value = -100
$ws.Cells.Item(14, 8).Copy($ws.Cells.Item(26, 12))
$ws.Cells.Item(26, 12).Value = -100
$ws.Cells.Item(14, 7).Copy($ws.Cells.Item(27, 3))
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(14, 7).Copy($ws.Cells.Item(27, 9))
$ws.Cells.Item(27, 9).Value = 1

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Cells.Item(16, 3).Value = 11
$ws.Cells.Item(16, 4).Value = 8
$ws.Cells.Item(16, 5).Value = 37.5
$ws.Cells.Item(16, 6).Value = 52
$ws.Cells.Item(16, 7).Value = 17
$ws.Cells.Item(16, 8).Value = 205.882352941176
$ws.Cells.Item(16, 9).Value = 25
$ws.Cells.Item(16, 10).Value = 10
$ws.Cells.Item(16, 11).Value = 150
$ws.Cells.Item(16, 12).Value = 92.307692307692
$ws.Cells.Item(16, 13).Value = 78.571428571428
$ws.Cells.Item(16, 14).Value = -60.9375
$ws.Cells.Item(17, 3).Value = 7
$ws.Cells.Item(17, 4).Value = 12
$ws.Cells.Item(17, 5).Value = -41.666666666666
$ws.Cells.Item(17, 6).Value = 39
$ws.Cells.Item(17, 7).Value = 58
$ws.Cells.Item(17, 8).Value = -32.758620689655
$ws.Cells.Item(17, 9).Value = 22
$ws.Cells.Item(17, 10).Value = 28
$ws.Cells.Item(17, 11).Value = -21.428571428571
$ws.Cells.Item(17, 12).Value = 83.333333333333
$ws.Cells.Item(17, 13).Value = 266.666666666667
$ws.Cells.Item(17, 14).Value = -8.333333333333
$ws.Cells.Item(18, 3).Value = 8
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 166.666666666667
$ws.Cells.Item(18, 6).Value = 26
$ws.Cells.Item(18, 7).Value = 11
$ws.Cells.Item(18, 8).Value = 136.363636363636
$ws.Cells.Item(18, 9).Value = 18
$ws.Cells.Item(18, 10).Value = 5
$ws.Cells.Item(18, 11).Value = 260
$ws.Cells.Item(18, 13).Value = 260
$ws.Cells.Item(18, 14).Value = -60
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 4).Value = 6
$ws.Cells.Item(19, 6).Value = 42
$ws.Cells.Item(19, 8).Value = 35.483870967741
$ws.Cells.Item(19, 9).Value = 16
$ws.Cells.Item(19, 10).Value = 16
$ws.Cells.Item(19, 12).Value = -20
$ws.Cells.Item(19, 13).Value = 33.333333333333
$ws.Cells.Item(19, 14).Value = -36
$ws.Cells.Item(20, 3).Value = 8
$ws.Cells.Item(20, 4).Value = 5
$ws.Cells.Item(20, 5).Value = 60
$ws.Cells.Item(20, 6).Value = 22
$ws.Cells.Item(20, 7).Value = 31
$ws.Cells.Item(20, 8).Value = -29.032258064516
$ws.Cells.Item(20, 9).Value = 13
$ws.Cells.Item(20, 10).Value = 8
$ws.Cells.Item(20, 11).Value = 62.5
$ws.Cells.Item(20, 12).Value = 44.444444444444
$ws.Cells.Item(20, 13).Value = 550
$ws.Cells.Item(20, 14).Value = -55.172413793103
$ws.Cells.Item(21, 3).Value = 40
$ws.Cells.Item(21, 5).Value = 17.647058823529
$ws.Cells.Item(21, 6).Value = 182
$ws.Cells.Item(21, 7).Value = 153
$ws.Cells.Item(21, 8).Value = 18.954248366013
$ws.Cells.Item(21, 9).Value = 94
$ws.Cells.Item(21, 10).Value = 68
$ws.Cells.Item(21, 11).Value = 38.235294117647
$ws.Cells.Item(21, 12).Value = 46.875
$ws.Cells.Item(21, 13).Value = 129.268292682927
$ws.Cells.Item(21, 14).Value = -51.295336787564
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(24, 3).Value = 18
$ws.Cells.Item(24, 4).Value = 21
$ws.Cells.Item(24, 5).Value = -14.285714285714
$ws.Cells.Item(24, 6).Value = 82
$ws.Cells.Item(24, 7).Value = 91
$ws.Cells.Item(24, 8).Value = -9.890109890109
$ws.Cells.Item(24, 9).Value = 34
$ws.Cells.Item(24, 10).Value = 32
$ws.Cells.Item(24, 11).Value = 6.25
$ws.Cells.Item(24, 12).Value = 30.769230769230
$ws.Cells.Item(24, 13).Value = 30.769230769230
$ws.Cells.Item(25, 3).Value = 18
$ws.Cells.Item(25, 4).Value = 11
$ws.Cells.Item(25, 5).Value = 63.636363636363
$ws.Cells.Item(25, 6).Value = 71
$ws.Cells.Item(25, 7).Value = 67
$ws.Cells.Item(25, 8).Value = 5.970149253731
$ws.Cells.Item(25, 9).Value = 38
$ws.Cells.Item(25, 10).Value = 28
$ws.Cells.Item(25, 11).Value = 35.714285714285
$ws.Cells.Item(25, 12).Value = 31.034482758620
$ws.Cells.Item(25, 13).Value = 46.153846153846
$ws.Cells.Item(27, 4).Value = 3
$ws.Cells.Item(27, 5).Value = -66.666666666666
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(27, 7).Value = 5
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 5
$ws.Cells.Item(27, 11).Value = -80
$ws.Cells.Item(27, 12).Value = -50
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 8).Value = -83.333333333333
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 8).Value = -66.666666666666
